$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.118.40'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.279.70'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '155.47'
$ws.Range('E5').Value = '  +15,425.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '305.23'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '94.84'
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.493'
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '35.60'
$ws.Range('E11').Value = '  +8.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0803'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.633.54'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.44'
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.274.85'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('E18').Value = '  +4.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.043.51'
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.82'
$ws.Range('E20').Value = '  +4.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0917'
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.10'
$ws.Range('E23').Value = '  +1.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '243.73'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.95'
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.43'
$ws.Range('E29').Value = '  +6.89%  '
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('E31').Value = '  +1.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.64'
$ws.Range('E32').Value = '  +1.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.36'
$ws.Range('E33').Value = '  +3.31%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0753'
$ws.Range('E35').Value = '  +0.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.08'
$ws.Range('E36').Value = '  +1.26%  '
$ws.Range('E37').Value = '  +3.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.13'
$ws.Range('E38').Value = '  +2.34%  '
$ws.Range('E39').Value = '  -0.29%  '
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.19'
$ws.Range('E42').Value = '  +6.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.019.04'
$ws.Range('E43').Value = '  -2.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.60'
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('E45').Value = '  +11.38%  '
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.22'
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.99'
$ws.Range('E48').Value = '  +2.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.59'
$ws.Range('E49').Value = '  +3.46%  '
$ws.Range('E50').Value = '  -0.54%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.36'
$ws.Range('E51').Value = '  -0.63%  '
